# Regenerate the "K" column (column G) values on the active worksheet.
# This mirrors an upstream data-regeneration run (new strike-count / K
# values recomputed from std/mean and written back over the prior
# Strike# values) by writing the freshly computed values directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,40,41,42,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64)
$values = @(0,0,2,0,1,0,0,1,1,0,2,1,1,2,0,1,2,1,1,2,2,2,0,2,3,0,0,1,0,0,1,0,0,2,2,0,0,1,3,1,0,2,1,1,2,1,0,3,1,2,0,1,0,2,0,1,1,1,1,1,1)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $k = $values[$i]
    # Column G is the 7th column ("K" header in row 1).
    $ws.Cells.Item($r, 7).Value = $k
}
